# Update existing category rows (Windbreaker/Varsity/Mesh -> Kemeja/Sleeve/Topi)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Kemeja"
$ws.Range("C3").Value = "kemeja"

$ws.Range("B4").Value = "Sleeve"
$ws.Range("C4").Value = "sleeve"

$ws.Range("B5").Value = "Topi"
$ws.Range("C5").Value = "topi"

# Add two new category rows
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Tas"
$ws.Range("C6").Value = "tas"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Aksesoris"
$ws.Range("C7").Value = "aksesoris"

# Update selection to match final state
$ws.Range("B2").Select() | Out-Null
